$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 4998.5713
$ws.Range("I28").Value = 1724.375
$ws.Range("J28").Value = 9364.166999999999
$ws.Range("K28").Value = 1724.375
$ws.Range("L28").Value = 9364.166999999999
$ws.Range("M28").Value = -1239.375
$ws.Range("N28").Value = -10334.167
$ws.Range("H34").Value = 7370
$ws.Range("I34").Value = 7370
$ws.Range("K34").Value = 7370
$ws.Range("M34").Value = -7167
$ws.Range("H36").Value = 7370
$ws.Range("I36").Value = 7370
$ws.Range("K36").Value = 7370
$ws.Range("M36").Value = -6655
$ws.Range("H43").Value = 1622.3334
$ws.Range("I43").Value = 983.5
$ws.Range("K43").Value = 983.5
$ws.Range("M43").Value = -914.5
$ws.Range("H92").Value = 260
$ws.Range("J92").Value = 325
$ws.Range("L92").Value = 325
$ws.Range("N92").Value = -2821
$ws.Range("H94").Value = 3975.4285
$ws.Range("I94").Value = 3975.4285
$ws.Range("K94").Value = 3975.4285
$ws.Range("M94").Value = -3524.4285
$ws.Range("H96").Value = 142857980
$ws.Range("I96").Value = 200000240
$ws.Range("K96").Value = 600000720
$ws.Range("M96").Value = -599999347
$ws.Range("I106").Value = 1500
$ws.Range("K106").Value = 1500
$ws.Range("M106").Value = -869
$ws.Range("H107").Value = 3166.5
$ws.Range("I107").Value = 3599.8
$ws.Range("K107").Value = 3599.8
$ws.Range("M107").Value = -1679.8
$ws.Range("H111").Value = 1685.9678
$ws.Range("J111").Value = 1849.963
$ws.Range("L111").Value = 5549.889
$ws.Range("N111").Value = -11683.889
$ws.Range("H132").Value = 4133.1724
$ws.Range("I132").Value = 1385.3043
$ws.Range("J132").Value = 14666.667
$ws.Range("K132").Value = 4155.9129
$ws.Range("L132").Value = 44000.001
$ws.Range("M132").Value = -1625.9129
$ws.Range("N132").Value = -49060.001
$ws.Range("H138").Value = 2526.0667
$ws.Range("I138").Value = 1384
$ws.Range("K138").Value = 4152
$ws.Range("M138").Value = 988

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H122").Value = 2990
$ws.Range("I122").Value = 2990
$ws.Range("K122").Value = 8970
$ws.Range("M122").Value = -6520

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 913
$ws.Range("I20").Value = 354.6
$ws.Range("J20").Value = 1471.4
$ws.Range("K20").Value = 354.6
$ws.Range("L20").Value = 1471.4
$ws.Range("M20").Value = -107.6
$ws.Range("N20").Value = -1965.4
$ws.Range("H54").Value = 4824
$ws.Range("I54").Value = 4824
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 4824
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4340
$ws.Range("N54").ClearContents()
$ws.Range("H94").Value = 395
$ws.Range("I94").Value = 395
$ws.Range("K94").Value = 395
$ws.Range("M94").Value = 56

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 790.5263
$ws.Range("I22").Value = 730.3077
$ws.Range("J22").Value = 921
$ws.Range("K22").Value = 730.3077
$ws.Range("L22").Value = 921
$ws.Range("M22").Value = -380.3077
$ws.Range("N22").Value = -1621
$ws.Range("H58").Value = 3224
$ws.Range("I58").Value = 1650.1818
$ws.Range("J58").Value = 8994.666999999999
$ws.Range("K58").Value = 1650.1818
$ws.Range("L58").Value = 8994.666999999999
$ws.Range("M58").Value = -1447.1818
$ws.Range("N58").Value = -9400.666999999999
$ws.Range("H99").Value = 3754.3333
$ws.Range("J99").Value = 3631.75
$ws.Range("L99").Value = 3631.75
$ws.Range("N99").Value = -6627.75
$ws.Range("H126").Value = 3754.3333
$ws.Range("J126").Value = 3631.75
$ws.Range("L126").Value = 10895.25
$ws.Range("N126").Value = -15835.25
$ws.Range("H136").Value = 3224
$ws.Range("I136").Value = 1650.1818
$ws.Range("J136").Value = 8994.666999999999
$ws.Range("K136").Value = 4950.5454
$ws.Range("L136").Value = 26984.001
$ws.Range("M136").Value = -2400.5454
$ws.Range("N136").Value = -32084.001

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 8774.5
$ws.Range("I9").Value = 10229.6
$ws.Range("J9").Value = 1499
$ws.Range("K9").Value = 30688.8
$ws.Range("L9").Value = 4497
$ws.Range("M9").Value = -30464.8
$ws.Range("N9").Value = -4945
$ws.Range("H63").Value = 1300
$ws.Range("I63").Value = 950
$ws.Range("J63").Value = 2000
$ws.Range("K63").Value = 2850
$ws.Range("L63").Value = 6000
$ws.Range("M63").Value = -2101
$ws.Range("N63").Value = -7498
$ws.Range("H66").Value = 1300
$ws.Range("I66").Value = 950
$ws.Range("J66").Value = 2000
$ws.Range("K66").Value = 8550
$ws.Range("L66").Value = 18000
$ws.Range("M66").Value = -4806
$ws.Range("N66").Value = -25488
$ws.Range("H86").Value = 558.9
$ws.Range("I86").Value = 554.3333
$ws.Range("K86").Value = 1662.9999
$ws.Range("M86").Value = -476.9999
$ws.Range("H89").Value = 558.9
$ws.Range("I89").Value = 554.3333
$ws.Range("K89").Value = 4988.9997
$ws.Range("M89").Value = 939.0002999999997
$ws.Range("H97").Value = 1875
$ws.Range("J97").Value = 1875
$ws.Range("L97").Value = 5625
$ws.Range("N97").Value = -6617
$ws.Range("H107").Value = 468.43478
$ws.Range("I107").Value = 231.36363
$ws.Range("J107").Value = 685.75
$ws.Range("K107").Value = 694.0908899999999
$ws.Range("L107").Value = 2057.25
$ws.Range("M107").Value = 1225.90911
$ws.Range("N107").Value = -5897.25
$ws.Range("H131").Value = 1924.3334
$ws.Range("I131").Value = 1387.25
$ws.Range("K131").Value = 4161.75
$ws.Range("M131").Value = 878.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1260
$ws.Range("I80").Value = 1350
$ws.Range("J80").Value = 1200
$ws.Range("K80").Value = 1350
$ws.Range("L80").Value = 1200
$ws.Range("M80").Value = -352
$ws.Range("N80").Value = -3196
$ws.Range("H83").Value = 1260
$ws.Range("I83").Value = 1350
$ws.Range("J83").Value = 1200
$ws.Range("K83").Value = 6750
$ws.Range("L83").Value = 6000
$ws.Range("M83").Value = -1758
$ws.Range("N83").Value = -15984
$ws.Range("H106").Value = 75000
$ws.Range("J106").Value = 75000
$ws.Range("L106").Value = 75000
$ws.Range("N106").Value = -77524
$ws.Range("H108").Value = 60684
$ws.Range("J108").Value = 60684
$ws.Range("L108").Value = 60684
$ws.Range("N108").Value = -68364
$ws.Range("H122").Value = 5000
$ws.Range("I122").Value = 4000
$ws.Range("K122").Value = 12000
$ws.Range("M122").Value = -9550

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5696.45
$ws.Range("I46").Value = 5499.857
$ws.Range("J46").Value = 5802.3076
$ws.Range("K46").Value = 5499.857
$ws.Range("L46").Value = 5802.3076
$ws.Range("M46").Value = -5311.857
$ws.Range("N46").Value = -6178.3076
$ws.Range("H55").Value = 919.0526
$ws.Range("I55").Value = 944.4
$ws.Range("J55").Value = 890.8889
$ws.Range("K55").Value = 944.4
$ws.Range("L55").Value = 890.8889
$ws.Range("M55").Value = -771.4
$ws.Range("N55").Value = -1236.8889
$ws.Range("H93").Value = 1390.1052
$ws.Range("I93").Value = 1501
$ws.Range("K93").Value = 1501
$ws.Range("M93").Value = -253

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 55556372
$ws.Range("I107").Value = 83333784
$ws.Range("J107").Value = 1545
$ws.Range("K107").Value = 250001352
$ws.Range("L107").Value = 4635
$ws.Range("M107").Value = -249999432
$ws.Range("N107").Value = -8475
$ws.Range("H113").Value = 709.1818
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 1080.2
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 3240.6
$ws.Range("M113").Value = 970
$ws.Range("N113").Value = -7580.6
$ws.Range("H122").Value = 1899.8462
$ws.Range("I122").Value = 1838.8889
$ws.Range("J122").Value = 2037
$ws.Range("K122").Value = 5516.6667
$ws.Range("L122").Value = 6111
$ws.Range("M122").Value = -3066.6667
$ws.Range("N122").Value = -11011
